$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.422602653503418
$ws.Range("B1").Value = 2.135307312011719
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.778645038604736
$ws.Range("E1").Value = 0.7235150337219238
